$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing maintenance-log entry ---
$ws.Cells.Item(2, 1).Value = "2025-03-26T12:37"
$ws.Cells.Item(2, 2).Value = "SHIFT_1"
$ws.Cells.Item(2, 3).Value = 0.95
$ws.Cells.Item(2, 4).Value = "OK"
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 8).Value = "Suriya"
$ws.Cells.Item(2, 9).Value = "Suriya"
$ws.Cells.Item(2, 10).Value = "Suriya"

# --- Row 3: new entry, same shift/readings as row 2 ---
$ws.Cells.Item(3, 1).Value = "2025-03-26T12:37"
$ws.Cells.Item(3, 2).Value = "SHIFT_1"
$ws.Cells.Item(3, 3).Value = 0.95
$ws.Cells.Item(3, 4).Value = "OK"
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = "Suriya"
$ws.Cells.Item(3, 9).Value = "Suriya"
$ws.Cells.Item(3, 10).Value = "Suriya"

# --- Row 4: another new entry, same shift/readings ---
$ws.Cells.Item(4, 1).Value = "2025-03-26T12:37"
$ws.Cells.Item(4, 2).Value = "SHIFT_1"
$ws.Cells.Item(4, 3).Value = 0.95
$ws.Cells.Item(4, 4).Value = "OK"
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = "Suriya"
$ws.Cells.Item(4, 9).Value = "Suriya"
$ws.Cells.Item(4, 10).Value = "Suriya"

# --- Row 5: a failing (NOK) entry. Numeric-looking values are entered as
# text (column formatted as Text before the write, then restored to the
# Normal style) so Excel does not auto-convert "0.95"/"1.2" into numbers. ---
$ws.Cells.Item(5, 1).Value = "2025-03-26T12:41"
$ws.Cells.Item(5, 2).Value = "SHIFT_1"

$c = $ws.Cells.Item(5, 3)
$c.NumberFormat = "@"
$c.Value = "0.95"
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "1.2"
$c.Style = "Normal"

$ws.Cells.Item(5, 5).Value = "NOK"
$ws.Cells.Item(5, 6).Value = "NOK"
$ws.Cells.Item(5, 7).Value = "NOK"
$ws.Cells.Item(5, 8).Value = "t"
$ws.Cells.Item(5, 9).Value = "t"
$ws.Cells.Item(5, 10).Value = "t"
